# Update the weekly price rows (3-7) for Fruta / Hortofrutícola Agro Chillán - Caqui
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44330
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("S3").Value = 861

# Row 4
$ws.Range("D4").Value = 44334
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("Q4").Value = "$/caja 12 kilos empedrada"
$ws.Range("S4").Value = 1042
$ws.Range("T4").Value = 12

# Row 5
$ws.Range("D5").Value = 44344
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 750
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 17500
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17750
$ws.Range("Q6").Value = "$/caja 16 kilos granel"
$ws.Range("S6").Value = 1109
$ws.Range("T6").Value = 16

# Row 7
$ws.Range("D7").Value = 44316
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 16000
$ws.Range("P7").Value = 16000
$ws.Range("Q7").Value = "$/caja 16 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 16
